$wb = $excel.ActiveWorkbook

# --- 1. Rename "Hoja1" -> "Aulas" ---
$wsAulas = $wb.Worksheets.Item("Hoja1")
$wsAulas.Name = "Aulas"

$wsProf = $wb.Worksheets.Item("Profesores")

# --- 2. Wipe existing data on both sheets ---
$wsAulas.Rows("1:3").Delete()
$wsProf.Rows("1:1").Delete()

# --- 3. "Aulas" sheet: 2 rows x 6 cols ---
# Numeric cells: assign Value first (stays a real number).
$wsAulas.Range("A1").Value = 8484
$wsAulas.Range("D1").Value = 30

# Text cells (including numeric-looking text): force Text format BEFORE
# assigning the value so it is stored as a shared string, not a number.
$wsAulas.Range("B1").NumberFormat = "@"
$wsAulas.Range("B1").Value = "A201-B"
$wsAulas.Range("C1").NumberFormat = "@"
$wsAulas.Range("C1").Value = "Laboratorio"
$wsAulas.Range("E1").NumberFormat = "@"
$wsAulas.Range("E1").Value = "idk"
$wsAulas.Range("F1").NumberFormat = "@"
$wsAulas.Range("F1").Value = "E201"

$wsAulas.Range("A2").NumberFormat = "@"
$wsAulas.Range("A2").Value = "7845"
$wsAulas.Range("B2").NumberFormat = "@"
$wsAulas.Range("B2").Value = "A202-C"
$wsAulas.Range("C2").NumberFormat = "@"
$wsAulas.Range("C2").Value = "Salon"
$wsAulas.Range("D2").NumberFormat = "@"
$wsAulas.Range("D2").Value = "30"
$wsAulas.Range("E2").NumberFormat = "@"
$wsAulas.Range("E2").Value = "ayuda"
$wsAulas.Range("F2").NumberFormat = "@"
$wsAulas.Range("F2").Value = "me"

# Make sure the whole block ends up Text-formatted (matches target style).
$wsAulas.Range("A1:F2").NumberFormat = "@"

# --- 4. "Profesores" sheet: 2 rows x 5 cols ---
$wsProf.Range("A1").Value = 1546
$wsProf.Range("B1").Value = 1

$wsProf.Range("C1").NumberFormat = "@"
$wsProf.Range("C1").Value = "Fulanito Martinez"
$wsProf.Range("D1").NumberFormat = "@"
$wsProf.Range("D1").Value = "ZXC"
$wsProf.Range("E1").NumberFormat = "@"
$wsProf.Range("E1").Value = "BNM"

$wsProf.Range("A2").Value = 6465

$wsProf.Range("B2").NumberFormat = "@"
$wsProf.Range("B2").Value = "1"
$wsProf.Range("C2").NumberFormat = "@"
$wsProf.Range("C2").Value = "Arturo Perez Reverte"
$wsProf.Range("D2").NumberFormat = "@"
$wsProf.Range("D2").Value = "QWE"
$wsProf.Range("E2").NumberFormat = "@"
$wsProf.Range("E2").Value = "RTY"

$wsProf.Range("A1:E2").NumberFormat = "@"

# Column widths: engine's ColumnWidth setter adds a fixed +5/6 padding
# before writing <col width=.../>, so back it out to land on the target
# bestFit-style widths (5, 2, 18, 5.109375, 5.109375).
$wsProf.Columns("A").ColumnWidth = 4.166666666666667
$wsProf.Columns("B").ColumnWidth = 1.166666666666667
$wsProf.Columns("C").ColumnWidth = 17.166666666666667
$wsProf.Columns("D").ColumnWidth = 4.25
$wsProf.Columns("E").ColumnWidth = 4.25

# --- 5. Add the new empty "Hoja2" sheet at the end ---
$wsHoja2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsProf)
$wsHoja2.Name = "Hoja2"

# --- 6. Selection / view cosmetics (best effort) ---
$wsAulas.Range("D12").Select()
$wsProf.Range("E2").Select()
$wsAulas.Select()
